$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.Value = "'56.409.21"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(2, 5)
$cell.Value = "'  +4.05%  "
$cell.Style = "Normal"

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.Value = "'2.974.01"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 5)
$cell.Value = "'  +3.50%  "
$cell.Style = "Normal"

# Row 4
$cell = $ws.Cells.Item(4, 4)
$cell.Value = "'1.00"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(4, 5)
$cell.Value = "'  -0.03%  "
$cell.Style = "Normal"

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.Value = "'501.75"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 5)
$cell.Value = "'  +5.62%  "
$cell.Style = "Normal"

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.Value = "'134.59"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 5)
$cell.Value = "'  +6.78%  "
$cell.Style = "Normal"

# Row 8
$cell = $ws.Cells.Item(8, 5)
$cell.Value = "'  +5.86%  "
$cell.Style = "Normal"

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.Value = "'7.44"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 5)
$cell.Value = "'  +11.58%  "
$cell.Style = "Normal"

# Row 10
$cell = $ws.Cells.Item(10, 5)
$cell.Value = "'  +9.54%  "
$cell.Style = "Normal"

# Row 11
$cell = $ws.Cells.Item(11, 5)
$cell.Value = "'  +4.40%  "
$cell.Style = "Normal"

# Row 12
$cell = $ws.Cells.Item(12, 5)
$cell.Value = "'  +3.30%  "
$cell.Style = "Normal"

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.Value = "'3.477.41"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 5)
$cell.Value = "'  +3.11%  "
$cell.Style = "Normal"

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.Value = "'25.34"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 5)
$cell.Value = "'  +11.11%  "
$cell.Style = "Normal"

# Row 15
$cell = $ws.Cells.Item(15, 5)
$cell.Value = "'  +11.22%  "
$cell.Style = "Normal"

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.Value = "'56.368.17"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 5)
$cell.Value = "'  +4.01%  "
$cell.Style = "Normal"

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.Value = "'2.972.04"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 5)
$cell.Value = "'  +2.95%  "
$cell.Style = "Normal"

# Row 18
$cell = $ws.Cells.Item(18, 4)
$cell.Value = "'5.73"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.Value = "'  +9.27%  "
$cell.Style = "Normal"

# Row 19
$cell = $ws.Cells.Item(19, 5)
$cell.Value = "'  +6.08%  "
$cell.Style = "Normal"

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.Value = "'7.72"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 5)
$cell.Value = "'  +8.75%  "
$cell.Style = "Normal"

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.Value = "'321.91"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 5)
$cell.Value = "'  +4.10%  "
$cell.Style = "Normal"

# Row 22
$cell = $ws.Cells.Item(22, 5)
$cell.Value = "'  +0.33%  "
$cell.Style = "Normal"

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.Value = "'0.468"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 5)
$cell.Value = "'  +4.27%  "
$cell.Style = "Normal"

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.Value = "'61.72"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.Value = "'  +3.34%  "
$cell.Style = "Normal"

# Row 25
$cell = $ws.Cells.Item(25, 5)
$cell.Value = "'  +0.06%  "
$cell.Style = "Normal"

# Row 26
$cell = $ws.Cells.Item(26, 5)
$cell.Value = "'  +5.46%  "
$cell.Style = "Normal"

# Row 27
$cell = $ws.Cells.Item(27, 5)
$cell.Value = "'  +7.64%  "
$cell.Style = "Normal"

# Row 28
$cell = $ws.Cells.Item(28, 5)
$cell.Value = "'  +1.63%  "
$cell.Style = "Normal"

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.Value = "'6.74"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 5)
$cell.Value = "'  +8.43%  "
$cell.Style = "Normal"

# Row 30
$cell = $ws.Cells.Item(30, 2)
$cell.Value = "'Fetch.AI"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 3)
$cell.Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 4)
$cell.Value = "'1.18"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 5)
$cell.Value = "'  +2.55%  "
$cell.Style = "Normal"

# Row 31
$cell = $ws.Cells.Item(31, 2)
$cell.Value = "'PancakeSwap"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 3)
$cell.Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 4)
$cell.Value = "'1.74"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 5)
$cell.Value = "'  +7.53%  "
$cell.Style = "Normal"

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.Value = "'20.31"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 5)
$cell.Value = "'  +6.35%  "
$cell.Style = "Normal"

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.Value = "'158.37"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 5)
$cell.Value = "'  +14.33%  "
$cell.Style = "Normal"

# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.Value = "'4.43"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 5)
$cell.Value = "'  +4.20%  "
$cell.Style = "Normal"

# Row 35
$cell = $ws.Cells.Item(35, 2)
$cell.Value = "'Aptos"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 3)
$cell.Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 4)
$cell.Value = "'5.52"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 5)
$cell.Value = "'  +0.99%  "
$cell.Style = "Normal"

# Row 36
$cell = $ws.Cells.Item(36, 2)
$cell.Value = "'ImmutableX"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 3)
$cell.Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 4)
$cell.Value = "'1.25"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 5)
$cell.Value = "'  +3.02%  "
$cell.Style = "Normal"

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.Value = "'0.0670"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 5)
$cell.Value = "'  +7.82%  "
$cell.Style = "Normal"

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.Value = "'22.85"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 5)
$cell.Value = "'  -0.66%  "
$cell.Style = "Normal"

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.Value = "'3.005.40"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 5)
$cell.Value = "'  +3.46%  "
$cell.Style = "Normal"

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.Value = "'1.00"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 5)
$cell.Value = "'  +0.00%  "
$cell.Style = "Normal"

# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.Value = "'36.14"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 5)
$cell.Value = "'  +2.29%  "
$cell.Style = "Normal"

# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.Value = "'0.638"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 5)
$cell.Value = "'  +6.52%  "
$cell.Style = "Normal"

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.Value = "'2.233.60"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.Value = "'  +8.56%  "
$cell.Style = "Normal"

# Row 44
$cell = $ws.Cells.Item(44, 5)
$cell.Value = "'  +5.38%  "
$cell.Style = "Normal"

# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.Value = "'0.979"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 5)
$cell.Value = "'  +2.06%  "
$cell.Style = "Normal"

# Row 46
$cell = $ws.Cells.Item(46, 4)
$cell.Value = "'3.55"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 5)
$cell.Value = "'  +3.39%  "
$cell.Style = "Normal"

# Row 47
$cell = $ws.Cells.Item(47, 5)
$cell.Value = "'  +18.69%  "
$cell.Style = "Normal"

# Row 48
$cell = $ws.Cells.Item(48, 5)
$cell.Value = "'  +9.64%  "
$cell.Style = "Normal"

# Row 49
$cell = $ws.Cells.Item(49, 5)
$cell.Value = "'  +6.98%  "
$cell.Style = "Normal"

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.Value = "'18.88"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 5)
$cell.Value = "'  +5.41%  "
$cell.Style = "Normal"

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.Value = "'0.0862"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 5)
$cell.Value = "'  +8.28%  "
$cell.Style = "Normal"
